$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45853, 720,     477.5,             720,     477.5,             1640, 283),
    @(45854, 695.625, 482,               695.625, 482,               1499, 291),
    @(45855, 708.5,   492.5,             708.5,   492.5,             1499, 278),
    @(45856, 723,     560,               723,     560,               1859, 275),
    @(45857, 716,     985,               716,     985,               2840, 275),
    @(45859, 714,     808.4649999999999, 714,     808.4649999999999, 2840, 294)
)

$startRow = 21
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
